$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 167, pushing the existing rows 167:190 down to 168:191,
# then populate it with the new weekly price record.
$ws.Rows("167:167").Insert()

$ws.Range("A167").Value = 10
$ws.Range("B167").Value = "Vega Modelo de Temuco"
$ws.Range("C167").Value = "La Araucanía"
$ws.Range("D167").Value = 45127
$ws.Range("E167").Value = 9
$ws.Range("F167").Value = 100112035
$ws.Range("G167").Value = "Bruselas (repollito)"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 180
$ws.Range("K167").Value = 24000
$ws.Range("L167").Value = 25000
$ws.Range("M167").Value = 24556
$ws.Range("N167").Value = '$/malla 15 kilos'
$ws.Range("O167").Value = "Provincia de Quillota"
$ws.Range("P167").Value = 1637
$ws.Range("Q167").Value = 15
$ws.Range("R167").Value = "Hortaliza"
